# Tool-Comparison.xlsx edit:
#  - Insert a new blank column before column A (shifts the whole table right
#    from A:D to B:E, including column widths, the Table1 ListObject and the
#    conditional-formatting range).
#  - Add a note "adsf" in C16.
#  - Leave the final selection on D14 (matches the saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything one column to the right by inserting a new column A.
$ws.Range("A:A").EntireColumn.Insert()

# The table (now living at B1:E12) needs its own range definition refreshed
# so the <table> part / autoFilter ref follow the shifted data.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B1:E12"))

# Likewise, re-anchor the conditional formatting (still on the old B2:D12
# reference after the column insert) onto the new C2:E12 data range.
$fc = $ws.Range("B2:D12").FormatConditions
$fc.Item(1).ModifyAppliesToRange($ws.Range("C2:E12"))

# New note cell.
$ws.Range("C16").Value = "adsf"

# Final selection left on D14.
$ws.Range("D14").Select() | Out-Null
